$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused rows 3-10 first (shifts nothing else, just shrinks used range)
$ws.Rows("3:10").Delete()

# --- Row 1 (headers) ---
$ws.Range("B1").Value = "terminalID"
$ws.Range("F1").Value = "dateCreated"
$ws.Range("G1").Value = "institutionID"
$ws.Range("H1").Value = "profileName"

# --- Row 2 (data) ---
$ws.Range("B2").Value = 30495868
$ws.Range("C2").Value = "topwise"
$ws.Range("D2").Value = 122345
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2
$ws.Range("G2").ClearFormats()
$ws.Range("G2").Value = "FREE380827"
$ws.Range("H2").Value = "EPMS"

# Update selection to match the new active cell
$ws.Range("C1").Select() | Out-Null
